$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: remove scope param, change treeId 7 -> 8
$ws.Range("D2").Value = "i=1&p=25&s%5BL3%5D=-1&f%5Bfamily%5D%5Bvalue%5D%5B%5D=RH_PRODUCTS_TRADEITEMS_VARIANTS&f%5Bfamily%5D%5Btype%5D=in&f%5Bcategory%5D%5Bvalue%5D%5BtreeId%5D=8&f%5Bcategory%5D%5Bvalue%5D%5BcategoryId%5D=0&f%5Bcategory%5D%5Btype%5D=1&t=product-grid"

# Row 3: remove __L3, __PARENTID, __SOURCEITEMID params; scope channel_product_catalog -> ecommerce; treeId 79 -> 88
$ws.Range("D3").Value = "i=1&p=25&s%5BL3%5D=-1&f%5Bfamily%5D%5Bvalue%5D%5B%5D=RH_PRODUCTS_TRADEITEMS_VARIANTS&f%5Bfamily%5D%5Btype%5D=in&f%5Bscope%5D%5Bvalue%5D=ecommerce&f%5Bcategory%5D%5Bvalue%5D%5BtreeId%5D=88&f%5Bcategory%5D%5Bvalue%5D%5BcategoryId%5D=0&f%5Bcategory%5D%5Btype%5D=1&t=product-grid"

# Row 4: add s[updated]=1; scope channel_gs1 -> GS1_GDSN; remove __PARENTID
$ws.Range("D4").Value = "i=1&p=25&s%5Bupdated%5D=1&f%5Bfamily%5D%5Bvalue%5D%5B%5D=GS1_GDSN&f%5Bfamily%5D%5Btype%5D=in&f%5Bscope%5D%5Bvalue%5D=GS1_GDSN&f%5Bcategory%5D%5Bvalue%5D%5BtreeId%5D=4&f%5Bcategory%5D%5Bvalue%5D%5BcategoryId%5D=0&f%5Bcategory%5D%5Btype%5D=1&t=product-grid"

# Row 5: scope channel_product_catalog -> ecommerce
$ws.Range("D5").Value = "i=1&p=25&s%5Bupdated%5D=1&f%5Bfamily%5D%5Bvalue%5D%5B%5D=MD_RECIPIENT_MAPPING&f%5Bfamily%5D%5Btype%5D=in&f%5Bscope%5D%5Bvalue%5D=ecommerce&f%5Bcategory%5D%5Bvalue%5D%5BtreeId%5D=1&f%5Bcategory%5D%5Bvalue%5D%5BcategoryId%5D=3&f%5Bcategory%5D%5Btype%5D=1&t=product-grid"

# Row 6: scope channel_product_catalog -> ecommerce
$ws.Range("D6").Value = "i=1&p=25&s%5Bupdated%5D=1&f%5Bfamily%5D%5Bvalue%5D%5B%5D=MD_SUPPLIER_MAPPING&f%5Bfamily%5D%5Btype%5D=in&f%5Bscope%5D%5Bvalue%5D=ecommerce&f%5Bcategory%5D%5Bvalue%5D%5BtreeId%5D=1&f%5Bcategory%5D%5Bvalue%5D%5BcategoryId%5D=3&f%5Bcategory%5D%5Btype%5D=1&t=product-grid"

# Row 7: i=3 -> i=1; scope channel_product_catalog -> ecommerce
$ws.Range("D7").Value = "i=1&p=25&s%5Bupdated%5D=1&f%5Bfamily%5D%5Bvalue%5D%5B%5D=MD_SUPPLIER_MASTER&f%5Bfamily%5D%5Btype%5D=in&f%5Bscope%5D%5Bvalue%5D=ecommerce&f%5Bcategory%5D%5Bvalue%5D%5BtreeId%5D=1&f%5Bcategory%5D%5Bvalue%5D%5BcategoryId%5D=3&f%5Bcategory%5D%5Btype%5D=1&t=product-grid"

# Row 8: s[updated]=1 -> s[MD_HUB_GLOBAL_ENTERPRISE_ID]=-1; scope channel_product_catalog -> PRODUCT_CATALOG
$ws.Range("D8").Value = "i=1&p=25&s%5BMD_HUB_GLOBAL_ENTERPRISE_ID%5D=-1&f%5Bfamily%5D%5Bvalue%5D%5B%5D=MD_HUB&f%5Bfamily%5D%5Btype%5D=in&f%5Bscope%5D%5Bvalue%5D=PRODUCT_CATALOG&f%5Bcategory%5D%5Bvalue%5D%5BtreeId%5D=1&f%5Bcategory%5D%5Bvalue%5D%5BcategoryId%5D=3&f%5Bcategory%5D%5Btype%5D=1&t=product-grid"

# New row 10: UNASSIGNED PRODUCTS AND ITEMS
$ws.Range("A10").Value = "UNASSIGNED PRODUCTS AND ITEMS"
$ws.Range("B10").Value = "updated,MD_SUPPLIER_MASTER_SUPPLIERID,MD_SUPPLIER_MASTER_SUPPLIERNAME,MD_SUPPLIER_MASTER_ALIASES,MD_SUPPLIER_MASTER_SUPPLIER_SITE,MD_SUPPLIER_MASTER_GLN"
$ws.Range("C10").Value = "product-grid"
$ws.Range("D10").Value = "i=1&p=25&s%5Bupdated%5D=1&f%5Bscope%5D%5Bvalue%5D=PRODUCT_CATALOG&f%5Bcategory%5D%5Bvalue%5D%5BtreeId%5D=1&f%5Bcategory%5D%5Bvalue%5D%5BcategoryId%5D=-1&f%5Bcategory%5D%5Btype%5D=1&t=product-grid"
$ws.Range("E10").Value = "admin"
$ws.Range("F10").Value = "public"
